$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.716657
$ws.Range("N2").Value = 5.149971
$ws.Range("O2").Value = 0.3840886036988016
$ws.Range("P2").Value = 0.3840886036988015
$ws.Range("Q2").Value = 0.224546746666
$ws.Range("R2").Value = 2.020920719994
$ws.Range("S2").Value = 0.3840886036988016
$ws.Range("T2").Value = 0.3840886036988015

# Row 3
$ws.Range("O3").Value = 0.07870146593648156
$ws.Range("P3").Value = 0.07870146593648154
$ws.Range("S3").Value = 0.07870146593648156
$ws.Range("T3").Value = 0.07870146593648154

# Row 4
$ws.Range("M4").Value = 1.677572333333333
$ws.Range("N4").Value = 5.032717
$ws.Range("O4").Value = 0.3753437146230962
$ws.Range("P4").Value = 0.3753437146230962
$ws.Range("Q4").Value = 0.2194342898708888
$ws.Range("R4").Value = 1.974908608838
$ws.Range("S4").Value = 0.3753437146230962
$ws.Range("T4").Value = 0.3753437146230962

# Row 5
$ws.Range("M5").Value = 0.7234496666666667
$ws.Range("N5").Value = 2.170349
$ws.Range("O5").Value = 0.1618662157416207
$ws.Range("P5").Value = 0.1618662157416207
$ws.Range("Q5").Value = 0.09463059249844444
$ws.Range("R5").Value = 0.8516753324859999
$ws.Range("S5").Value = 0.1618662157416207
$ws.Range("T5").Value = 0.1618662157416207
